$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.612.33"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.959.97"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "243.45"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "60.42"
$ws.Range("E7").Value = "  +5.57%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "0.375"
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("D10").Value = "0.0788"
$ws.Range("E10").Value = "  -6.78%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("E12").Value = "  +5.60%  "
$ws.Range("D13").Value = "21.81"
$ws.Range("E13").Value = "  +2.85%  "
$ws.Range("D14").Value = "2.246.81"
$ws.Range("E14").Value = "  +0.88%  "
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "1.959.36"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "36.507.37"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "69.62"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("D21").Value = "229.45"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "5.08"
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +3.71%  "
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("E26").Value = "  +5.30%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "161.26"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  +21.52%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "4.75"
$ws.Range("E32").Value = "  +4.04%  "
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").Value = "4.44"
$ws.Range("E34").Value = "  +6.48%  "
$ws.Range("D35").Value = "3.48"
$ws.Range("E35").Value = "  +9.97%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +4.61%  "
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("E39").Value = "  -12.16%  "
$ws.Range("D40").Value = "0.0967"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").Value = "1.17"
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "15.83"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").Value = "1.362.88"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").Value = "88.76"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("E47").Value = "  +0.10%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.10"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("D50").Value = "45.84"
$ws.Range("E50").Value = "  +6.15%  "
$ws.Range("D51").Value = "2.138.10"
$ws.Range("E51").Value = "  +0.95%  "
